$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of each data row (2-16) to the row whose D/J/K/L/M/P values it
# should now hold (a permutation of the original rows' tuples).
# newRow -> sourceRow(before the edit)
$rowData = @{
    2  = @{ D = 44504; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 }
    3  = @{ D = 44301; J = 40; K = 3000; L = 3000; M = 3000; P = 3000 }
    4  = @{ D = 44497; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
    5  = @{ D = 44509; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
    6  = @{ D = 44315; J = 40; K = 4000; L = 4000; M = 4000; P = 4000 }
    7  = @{ D = 44312; J = 50; K = 4000; L = 4000; M = 4000; P = 4000 }
    8  = @{ D = 44316; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
    9  = @{ D = 44280; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 }
    10 = @{ D = 44176; J = 10; K = 4000; L = 4000; M = 4000; P = 4000 }
    11 = @{ D = 44508; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 }
    12 = @{ D = 44291; J = 35; K = 4000; L = 4000; M = 4000; P = 4000 }
    13 = @{ D = 44498; J = 40; K = 4000; L = 4000; M = 4000; P = 4000 }
    14 = @{ D = 44259; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 }
    15 = @{ D = 44365; J = 55; K = 5000; L = 5000; M = 5000; P = 5000 }
    16 = @{ D = 44313; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J  # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K  # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals.L  # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals.M  # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals.P  # P: Precio $/Kg
}
